# Generate Report for Handoff
#
# The 4503c4c5-... file moved from "Handed back: in sync with en-US" to
# "Ready for handoff" status, and the handoff timestamps for both
# files were refreshed. Update the Overview sheet plus the per-language
# (zh-cn / de-de) detail sheets to match.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
# Row 3 corresponds to 4503c4c5-86d4-4d3b-9b35-781577df6db2.md -> status flips
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
# Latest Handoff Date refreshed for both rows
$ov.Range("D2").Value = "2016-51-14 09:51:56"
$ov.Range("D3").Value = "2016-51-14 09:51:56"

# --- zh-cn detail sheet -------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("E2").Value = "2016-03-14 09:51:49"
$zh.Range("E3").Value = "2016-03-14 09:51:49"

# --- de-de detail sheet -------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("E2").Value = "2016-03-14 09:51:56"
$de.Range("E3").Value = "2016-03-14 09:51:56"
